$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# ---- Row 2 ----
$ws.Range("A2").Value = "Gioavni"
$ws.Range("B2").Value = "'0756"
$ws.Range("C2").Value = "Escola Manoel Correia"
$ws.Range("D2").Value = "Sem comunicação de câmeras e tentar mudar para o DDNS."
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = "Pendente"
$ws.Range("H2").Value = "Maxvel: 21 / Forte: 14"
$ws.Rows.Item(2).AutoFit()

# ---- Row 3 ----
$ws.Range("A3").Value = "Giovani"
$ws.Range("B3").Value = "'0884"
$ws.Range("C3").Value = "Residencial Santos Drumont"
$ws.Range("D3").Value = "Sem comunicação geral, pode ser internet. Não consegui saber junto do cliente."
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = $null
$ws.Range("G3").Value = "Pendente"
$ws.Range("H3").Value = $null

# ---- Row 4 ----
$ws.Range("A4").Value = "Giovani"
$ws.Range("B4").Value = "'0840"
$ws.Range("C4").Value = "Valdemar Amaral"
$ws.Range("D4").Value = "Sem comunicação de alarmes, ontem a sala estava trancada e quem tinha a chave não estava lá na hora."
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("G4").Value = "Pendente"
$ws.Range("H4").Value = $null
$ws.Rows.Item(4).RowHeight = 30

# ---- Row 5 ----
$ws.Range("A5").Value = "Giovani"
$ws.Range("B5").Value = "'0079"
$ws.Range("C5").Value = "Med Center"
$ws.Range("D5").Value = "Sem comunicação de alarmes, Roberto esteve no local mas não conseguiu acesso a central."
$ws.Range("E5").Value = $null
$ws.Range("F5").Value = $null
$ws.Range("G5").Value = "Pendente"
$ws.Range("H5").Value = $null
$ws.Rows.Item(5).RowHeight = 30

# ---- Row 6 ----
$ws.Range("A6").Value = "Giovani"
$ws.Range("B6").Value = "'0288"
$ws.Range("C6").Value = "Viaceu Loja"
$ws.Range("D6").Value = "Sem comunicação de câmeras, cliente não responde e sem o minimo interesse em resolver."
$ws.Range("E6").Value = $null
$ws.Range("F6").Value = $null
$ws.Range("G6").Value = "Pendente"
$ws.Range("H6").Value = $null
$ws.Rows.Item(6).RowHeight = 30

# ---- Row 7 (keeps ht=30, do not touch row height) ----
$ws.Range("A7").Value = "Roberto"
$ws.Range("B7").Value = "'0390"
$ws.Range("C7").Value = "Igreja Imaculada"
$ws.Range("D7").Value = "Sem comunicação de alarmes, tentar atualizar a central pra gente poder ver se aceita módulo."
$ws.Range("E7").Value = $null
$ws.Range("F7").Value = $null
$ws.Range("G7").Value = "Pendente"
$ws.Range("H7").Value = $null

# ---- Row 8 ----
$ws.Range("A8").Value = "Roberto"
$ws.Range("B8").Value = "'0463"
$ws.Range("C8").Value = "Bc Refratário"
$ws.Range("D8").Value = "Zona 14 segue aberta e passar as câmeras para o DDNS."
$ws.Range("E8").Value = $null
$ws.Range("F8").Value = $null
$ws.Range("G8").Value = "Pendente"
$ws.Range("H8").Value = $null
$ws.Rows.Item(8).AutoFit()

# ---- Row 9 ----
$ws.Range("A9").Value = "Roberto"
$ws.Range("B9").Value = "'0422"
$ws.Range("C9").Value = "Sitio Alves"
$ws.Range("D9").Value = "Sem comunicação de câmeras e central via GPRS. Roberto não foi ontem."
$ws.Range("E9").Value = $null
$ws.Range("F9").Value = $null
$ws.Range("G9").Value = "Pendente"
$ws.Range("H9").Value = $null
$ws.Rows.Item(9).AutoFit()

# ---- Row 10 ----
$ws.Range("A10").Value = "Roberto"
$ws.Range("B10").Value = "'0554"
$ws.Range("C10").Value = "Valinhos Departamento Pessoal"
$ws.Range("D10").Value = "Sem comunicação de alarmes, cliente pediu reparo ontem. Linha telefônica."
$ws.Range("E10").Value = $null
$ws.Range("F10").Value = $null
$ws.Range("G10").Value = "Pendente"
$ws.Range("H10").Value = $null

# ---- Row 11 ----
$ws.Range("A11").Value = "Roberto"
$ws.Range("B11").Value = "'0217"
$ws.Range("C11").Value = "Brapi"
$ws.Range("D11").Value = "Zonas abertas, cliente pedindo reparo."
$ws.Range("E11").Value = $null
$ws.Range("F11").Value = $null
$ws.Range("G11").Value = $null
$ws.Range("H11").Value = $null

# ---- Row 12 (cleared entirely) ----
$ws.Range("A12").Value = $null
$ws.Range("B12").Value = $null
$ws.Range("C12").Value = $null
$ws.Range("D12").Value = $null
$ws.Range("E12").Value = $null
$ws.Range("F12").Value = $null
$ws.Range("G12").Value = $null
$ws.Range("H12").Value = $null

# ---- Selection ----
$ws.Range("H2").Select()

# ---- Window geometry (best effort; mirrors the saved workbookView) ----
$excel.ActiveWindow.Left = -120
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 20730
$excel.ActiveWindow.Height = 11040
